# Feature: Add joint borrowers support to Excel import with enterprise test data
#
# The sheet currently has 5 data rows (rows 2-5) and guarantor columns
# I..W (guarantors 1..5, 3 cols each). This edit:
#   1. Inserts 6 new columns before column I to hold two "joint borrower"
#      blocks (name / id / contact), pushing the existing guarantor
#      columns from I:W to O:AC.
#   2. Populates the two new header triplets (I1:N1).
#   3. Fills in the joint-borrower data for the existing rows (2-5) and
#      updates a couple of borrower records with new enterprise data.
#   4. Appends two brand new rows (6 and 7) with full data, including
#      enterprise (corporate) borrowers/joint-borrowers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 6 blank columns at I:N, shifting old I:W to O:AC -------
$ws.Range("I1:N1").EntireColumn.Insert()

# --- 2. New header cells for the two joint-borrower blocks -------------
$ws.Range("I1").Value = "共同借款人1名称"
$ws.Range("J1").Value = "共同借款人1证件号"
$ws.Range("K1").Value = "共同借款人1联系方式"
$ws.Range("L1").Value = "共同借款人2名称"
$ws.Range("M1").Value = "共同借款人2证件号"
$ws.Range("N1").Value = "共同借款人2联系方式"

# --- 3. Row 2 (城东支行 / 张三): add joint borrower 1 ------------------
$ws.Range("I2").Value = "赵强"
$ws.Range("J2").Value = "330107198406061111"
$ws.Range("K2").Value = "13700137011"

# --- Row 3 (滨江支行 / 李四): add joint borrower 2 ---------------------
$ws.Range("L3").Value = "黄丽"
$ws.Range("M3").Value = "330110198509092222"
$ws.Range("N3").Value = "13600136011"

# --- Row 4 (西湖支行): borrower becomes an enterprise, add joint
#     borrower 1 (also an enterprise) ----------------------------------
$ws.Range("B4").Value = "王五科技有限公司"
$ws.Range("C4").Value = "91330100MA2ABCDE12"
$ws.Range("D4").Value = "0571-88888888"
$ws.Range("E4").Value = "杭州市西湖区科技园A座"
$ws.Range("I4").Value = "杭州云端科技"
$ws.Range("J4").Value = "91330100MA2KLMNO56"
$ws.Range("K4").Value = "0571-77777777"

# --- Row 5 (萧山支行): joint borrower 1 replaces what used to be in the
#     first guarantor slot; the old guarantor 1 (吴刚) now correctly
#     sits in the (shifted) guarantor-1 columns (O5:Q5) ----------------
$ws.Range("I5").Value = "孙伟"
$ws.Range("J5").Value = "330108199211112222"
$ws.Range("K5").Value = "13700137012"

# --- 4. Brand new row 6: 拱墅支行 / 钱七实业集团 (enterprise) ----------
$ws.Range("A6").Value = "拱墅支行"
$ws.Range("B6").Value = "钱七实业集团"
$ws.Range("C6").Value = "91330100MA2FGHIJ34"
$ws.Range("D6").Value = "0571-99999999"
$ws.Range("E6").Value = "杭州市拱墅区工业路100号"
$ws.Range("I6").Value = "李明"
$ws.Range("J6").Value = "330109198512123333"
$ws.Range("K6").Value = "13700137013"
$ws.Range("L6").Value = "上海贸易有限公司"
$ws.Range("M6").Value = "91310000MA2PQRST78"
$ws.Range("N6").Value = "021-55555555"
$ws.Range("O6").Value = "郑伟"
$ws.Range("P6").Value = "330110198312125555"
$ws.Range("Q6").Value = "13700137005"
$ws.Range("R6").Value = "杭州融资担保公司"
$ws.Range("S6").Value = "91330100MA2UVWXY90"
$ws.Range("T6").Value = "0571-66666666"

# --- Brand new row 7: 余杭支行 / 孙八 -----------------------------------
$ws.Range("A7").Value = "余杭支行"
$ws.Range("B7").Value = "孙八"
$ws.Range("C7").Value = "330106198808083456"
$ws.Range("D7").Value = "13800138006"
$ws.Range("E7").Value = "杭州市余杭区BB街BB号"
$ws.Range("F7").Value = "周婷"
$ws.Range("G7").Value = "330106198909094567"
$ws.Range("H7").Value = "13900139006"
$ws.Range("L7").Value = "吴洁"
$ws.Range("M7").Value = "330111199012123333"
$ws.Range("N7").Value = "13600136013"
$ws.Range("O7").Value = "徐敏"
$ws.Range("P7").Value = "330111198411116666"
$ws.Range("Q7").Value = "13700137006"
$ws.Range("R7").Value = "马强"
$ws.Range("S7").Value = "330112198601017777"
$ws.Range("T7").Value = "13600136006"
